$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 so all existing data rows shift down by one
# (row 2 becomes row 3, row 3 becomes row 4, ... row 14 becomes row 15).
$ws.Rows.Item(2).Insert()

# The inserted row picks up the header row's formatting by default; strip it
# back to the plain (unstyled) look used by the other data rows.
$ws.Rows.Item(2).ClearFormats()

# Populate the newly inserted row with the latest date and the (constant)
# price values. Force the date cell to stay plain text, matching the other
# date cells in column A (which are stored as text, not real dates), then
# reset the cell style back to Normal so no residual number-format style
# lingers on the cell.
$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "2025-12-04"
$ws.Cells.Item(2, 1).Style = "Normal"

$ws.Cells.Item(2, 2).Value = 783.5
$ws.Cells.Item(2, 3).Value = 1112
$ws.Cells.Item(2, 4).Value = 3610
